# Apply updated odds values to Sheet1 per the diff (Jogos_da_Semana_FlashScore update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q7").Value = 1.77
$ws.Range("R7").Value = 1.97
$ws.Range("AR7").Value = 2.36
$ws.Range("G8").Value = 1.53
$ws.Range("H8").Value = 3.7
$ws.Range("J8").Value = 2.2
$ws.Range("L8").Value = 6.5
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("AA8").Value = 9
$ws.Range("AE8").Value = 8
$ws.Range("AP8").Value = 1.56
$ws.Range("AQ8").Value = 2.38
$ws.Range("AR8").Value = 2.95
$ws.Range("AS8").Value = 1.38
$ws.Range("G9").Value = 1.22
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("W9").Value = 2.5
$ws.Range("X9").Value = 1.5
$ws.Range("Y9").Value = 6
$ws.Range("AF9").Value = 10
$ws.Range("AM9").Value = 151
$ws.Range("H10").Value = 3.45
$ws.Range("I10").Value = 4.35
$ws.Range("J10").Value = 2.37
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 4.8
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 2.65
$ws.Range("Q10").Value = 2.05
$ws.Range("U10").Value = 1.45
$ws.Range("V10").Value = 2.37
$ws.Range("AE10").Value = 8.25
$ws.Range("AF10").Value = 6.8
$ws.Range("AM10").Value = 70
$ws.Range("G12").Value = 2.38
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3.2
$ws.Range("L12").Value = 3.75
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("Q12").Value = 2.3
$ws.Range("R12").Value = 1.6
$ws.Range("S12").Value = 4.33
$ws.Range("T12").Value = 1.2
$ws.Range("W12").Value = 1.91
$ws.Range("X12").Value = 1.8
$ws.Range("Y12").Value = 7
$ws.Range("Z12").Value = 11
$ws.Range("AA12").Value = 10
$ws.Range("AB12").Value = 23
$ws.Range("AC12").Value = 21
$ws.Range("AE12").Value = 7
$ws.Range("AI12").Value = 351
$ws.Range("AJ12").Value = 8.5
$ws.Range("AK12").Value = 15
$ws.Range("AL12").Value = 12
$ws.Range("AM12").Value = 34
$ws.Range("AN12").Value = 29
$ws.Range("AO12").Value = 41
$ws.Range("R18").Value = 1.33

Write-Output "Applied odds updates to rows 7, 8, 9, 10, 12, 18 on Sheet1."
